# Replace the literal "<br/>" markers in these six cells with real
# line breaks (Chr(10)), matching the author's re-authoring of the
# shared-string table (Small/Medium/Large Employees & Assets/Turnover
# criteria on the Thailand Summary sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("B24").Value = "<50 Production" + $nl + "<50 Services" + $nl + "<25 Wholesale" + $nl + "<15 Retail"
$ws.Range("C24").Value = "<50Millionlion bath Production, " + $nl + "<50Millionlion bath Services, " + $nl + "<50Millionlion bath Wholesale, " + $nl + "<30Millionlion bath Retail"

$ws.Range("B25").Value = "≤51 Production<200, " + $nl + "≤51 Services<200, " + $nl + "≤26 Wholesale<200, " + $nl + "≤16 Retail <30"
$ws.Range("C25").Value = "≤50 Production<200Millionlion bath, " + $nl + "≤50 Services<200Millionlion bath, " + $nl + "≤50 Wholesale<100Millionlion bath, " + $nl + "≤30 Retail <60Millionlion bath"

$ws.Range("B26").Value = ">=200 Production, " + $nl + ">=200 Services, " + $nl + ">=200 Wholesale, " + $nl + ">=30 Retail"
$ws.Range("C26").Value = ">=200Millionlion bath Production, " + $nl + ">=200Millionlion bath Services, " + $nl + ">=100Millionlion bath Wholesale, " + $nl + ">=60Millionlion bath Retail"
